$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - update column F (想去人数) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 770
$wsExhibit.Range("F3").Value = 0
$wsExhibit.Range("F4").Value = 0
$wsExhibit.Range("F5").Value = 0
$wsExhibit.Range("F6").Value = 0
$wsExhibit.Range("F7").Value = 0
$wsExhibit.Range("F8").Value = 0
$wsExhibit.Range("F9").Value = 0
$wsExhibit.Range("F11").Value = 530
$wsExhibit.Range("F12").Value = 61

# Sheet "全部类型" (sheet4.xml) - update column F (想去人数) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 770
$wsAll.Range("F3").Value = 0
$wsAll.Range("F7").Value = 4257
$wsAll.Range("F8").Value = 0
$wsAll.Range("F9").Value = 0
$wsAll.Range("F10").Value = 8613
$wsAll.Range("F11").Value = 222
$wsAll.Range("F12").Value = 0
$wsAll.Range("F16").Value = 0
$wsAll.Range("F17").Value = 0
